$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Row 19 (2016-11-21 / serial 42608): previously placeholder data (1, 0)
# gets replaced with the real observed totals.
# ---------------------------------------------------------------------
$ws.Range("B19").Value = 25
$ws.Range("C19").Value = 7

# ---------------------------------------------------------------------
# Row 20 (serial 42609): brand-new data row - fill in the raw values and
# copy the same formula pattern used by the rows above it.
# ---------------------------------------------------------------------
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 1

$ws.Range("D20").Formula = '=C20/B20'
$ws.Range("D20").NumberFormat = $ws.Range("D19").NumberFormat

$ws.Range("F20").Formula = '=(B19+B20)/2'
$ws.Range("G20").Formula = '=(B18+B20)/2'
$ws.Range("H20").Formula = '=(B17+B20)/2'
$ws.Range("I20").Formula = '=(B16+B20)/2'
$ws.Range("J20").Formula = '=(B15+B20)/2'
$ws.Range("K20").Formula = '=(B14+B20)/2'
$ws.Range("L20").Formula = '=(B13+B20)/2'

$ws.Range("M20").NumberFormat = $ws.Range("M19").NumberFormat

$ws.Range("N20").Formula = '=ABS($B20-F20)'
$ws.Range("O20").Formula = '=ABS($B20-G20)'
$ws.Range("P20").Formula = '=ABS($B20-H20)'
$ws.Range("Q20").Formula = '=ABS($B20-I20)'
$ws.Range("R20").Formula = '=ABS($B20-J20)'
$ws.Range("S20").Formula = '=ABS($B20-K20)'
$ws.Range("T20").Formula = '=ABS($B20-L20)'

# ---------------------------------------------------------------------
# Row 21 (serial 42610): another brand-new data row, same treatment.
# ---------------------------------------------------------------------
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 0

$ws.Range("D21").Formula = '=C21/B21'
$ws.Range("D21").NumberFormat = $ws.Range("D19").NumberFormat

$ws.Range("F21").Formula = '=(B20+B21)/2'
$ws.Range("G21").Formula = '=(B19+B21)/2'
$ws.Range("H21").Formula = '=(B18+B21)/2'
$ws.Range("I21").Formula = '=(B17+B21)/2'
$ws.Range("J21").Formula = '=(B16+B21)/2'
$ws.Range("K21").Formula = '=(B15+B21)/2'
$ws.Range("L21").Formula = '=(B14+B21)/2'

$ws.Range("M21").NumberFormat = $ws.Range("M19").NumberFormat

$ws.Range("N21").Formula = '=ABS($B21-F21)'
$ws.Range("O21").Formula = '=ABS($B21-G21)'
$ws.Range("P21").Formula = '=ABS($B21-H21)'
$ws.Range("Q21").Formula = '=ABS($B21-I21)'
$ws.Range("R21").Formula = '=ABS($B21-J21)'
$ws.Range("S21").Formula = '=ABS($B21-K21)'
$ws.Range("T21").Formula = '=ABS($B21-L21)'

# ---------------------------------------------------------------------
# Move the active selection back to A1 (top-left of the frozen pane).
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
